# Applies the cell-value updates for the crypto listing sheet (cryptos.xlsx).
#
# Cells in column D hold price text that can look like a plain number
# (e.g. "1.00", "0.0000200", "443.44"). If such a string is assigned directly
# via Range.Value, Excel auto-converts it to a numeric cell (stripping trailing
# zeros / turning it into scientific notation), which would not match the
# original text. To keep those cells as text (matching the source data, which
# stores every cell as a string), we momentarily force a text number format
# before assigning the value, then clear that formatting again afterwards so
# the cell's style index is left exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.638.76'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '3.109.94'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.22'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '615.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('E7').Value = '  -4.14%  '
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.106.35'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.739'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.203'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('E14').Value = '  +2.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.47'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').Value = '91.379.19'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '3.108.50'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.71'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.76'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '443.44'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.25'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000200'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -7.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.62'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.73'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.63'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.23%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +24.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.229'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.88%  '
$ws.Range('E32').Value = '  -9.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.177'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.29'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.981'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.68'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.13'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.07'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.93'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '488.97'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.435'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.39'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -6.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.21'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.73'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.14'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.56%  '
